$wb = $excel.ActiveWorkbook

# --- "ODI Batting": clear the blank placeholder INNING_NUMBER cells (B2, B3) ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B3").ClearContents()

# --- Add the new "ODI Batting Extra" sheet, placed after "ODI Bowling" ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$newSheet = $wb.Worksheets.Add($null, $odiBowling)
$newSheet.Name = "ODI Batting Extra"

# Reuse the bold/centered/bordered header style already used by the other
# sheets' row 1 (copy formats only, so the existing style index is shared
# instead of a new one being created).
$odiBowling.Range("A1:G1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Data rows: MATCH_CODE (A) is textual, BATTING_POSITION (B) numeric,
# NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL (C/D/E) textual (blank for most rows),
# MAN_OF_MATCH (F) textual "NO".
$dataRange = $newSheet.Range("A2:F4")
$dataRange.NumberFormat = "@"

$newSheet.Range("A2").Value = "4657"
$newSheet.Range("B2").Value = 8
$newSheet.Range("C2").Value = ""
$newSheet.Range("D2").Value = ""
$newSheet.Range("E2").Value = ""
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4658"
$newSheet.Range("B3").Value = 8
$newSheet.Range("C3").Value = ""
$newSheet.Range("D3").Value = ""
$newSheet.Range("E3").Value = ""
$newSheet.Range("F3").Value = "NO"

$newSheet.Range("A4").Value = "4679"
$newSheet.Range("B4").Value = 7
$newSheet.Range("C4").Value = "0"
$newSheet.Range("D4").Value = "0"
$newSheet.Range("E4").Value = ""
$newSheet.Range("F4").Value = "NO"

# BATTING_POSITION stays numeric.
$newSheet.Range("B2:B4").NumberFormat = "General"
$newSheet.Range("B2").Value = 8
$newSheet.Range("B3").Value = 8
$newSheet.Range("B4").Value = 7
